# Applies the update described by the diff: row 22 (T021), row 23 (T022),
# and row 30 (T029) of the "tareas" sheet get revised status/progress/date
# values, plus a refreshed description for T021.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tareas")

# Row 22 (T021): new description, status back to "En curso", 50% avance,
# fecha inicio moved one day later (2025-04-27).
$ws.Range("D22").Value = "Vacios plta mqta-protocolo de prueba hermeticidad hidrociclon 1-mqta"
$ws.Range("E22").Value = "En curso"
$ws.Range("F22").Value = 50
$ws.Range("H22").Value = Get-Date -Year 2025 -Month 4 -Day 27 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

# Row 23 (T022): status back to "Pausada", 40% avance.
$ws.Range("E23").Value = "Pausada"
$ws.Range("F23").Value = 40

# Row 30 (T029): avance raised to 80%.
$ws.Range("F30").Value = 80

# Restore the view position/selection captured in the saved workbook.
$excel.ActiveWindow.ScrollRow = 13
$ws.Range("E23").Select()
